$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.039.25"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.872.97"
$ws.Range("E3").Value = "  -2.61%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'319.55"
$ws.Range("E5").Value = "  -3.41%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.5056"
$ws.Range("E7").Value = "  -3.23%  "
$ws.Range("D8").Value = "'0.3942"
$ws.Range("E8").Value = "  -4.04%  "
$ws.Range("D9").Value = "'0.08210"
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").Value = "'42.21"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").Value = "'1.094"
$ws.Range("E11").Value = "  -3.20%  "
$ws.Range("D12").Value = "'23.84"
$ws.Range("E12").Value = "  +5.97%  "
$ws.Range("D13").Value = "1.866.99"
$ws.Range("E13").Value = "  -2.85%  "
$ws.Range("D14").Value = "'6.305"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").Value = "'7.190"
$ws.Range("E15").Value = "  -3.19%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'92.13"
$ws.Range("E17").Value = "  -4.45%  "
$ws.Range("D18").Value = "'0.00001091"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").Value = "'0.06424"
$ws.Range("E19").Value = "  -4.19%  "
$ws.Range("D20").Value = "'18.13"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D22").Value = "30.032.59"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -3.75%  "
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").Value = "'2.172"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").Value = "2.087.16"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").Value = "'21.38"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "'160.21"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").Value = "'2.234"
$ws.Range("E29").Value = "  -9.46%  "
$ws.Range("D30").Value = "'127.51"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").Value = "'1.065"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").Value = "'0.1035"
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("D33").Value = "'5.954"
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("D34").Value = "'3.687"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'0.02430"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("D36").Value = "'5.242"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "'0.06381"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("D38").Value = "'0.2149"
$ws.Range("E38").Value = "  -3.43%  "
$ws.Range("D39").Value = "'1.177"
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("D40").Value = "'8.510"
$ws.Range("E40").Value = "  -5.42%  "
$ws.Range("D41").Value = "'0.6317"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").Value = "'1.219"
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.03"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5925"
$ws.Range("E46").Value = "  -4.18%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.080"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.633"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("D49").Value = "'123.12"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'1.121"
$ws.Range("E51").Value = "  -3.38%  "
